$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns remain text (avoid Excel auto-numeric conversion)
$ws.Range("B2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "29.413.20"
$ws.Range("E2").Value = "  +0.05%  "
$ws.Range("D3").Value = "1.851.61"
$ws.Range("E3").Value = "  +0.25%  "
$ws.Range("D4").Value = "0.9995"
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "241.04"
$ws.Range("E5").Value = "  +0.18%  "
$ws.Range("D6").Value = "0.6292"
$ws.Range("D7").Value = "1.000"
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "0.07698"
$ws.Range("E8").Value = "  +1.91%  "
$ws.Range("D9").Value = "0.2937"
$ws.Range("E9").Value = "  -0.56%  "
$ws.Range("D10").Value = "24.58"
$ws.Range("E10").Value = "  +0.43%  "
$ws.Range("D11").Value = "0.07752"
$ws.Range("E11").Value = "  +0.78%  "
$ws.Range("D12").Value = "1.855.33"
$ws.Range("E12").Value = "  +0.09%  "
$ws.Range("B13").Value = "ShibaInu"
$ws.Range("C13").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D13").Value = "0.00001107"
$ws.Range("E13").Value = "  +10.30%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "5.026"
$ws.Range("E14").Value = "  +0.80%  "
$ws.Range("E15").Value = "  -0.37%  "
$ws.Range("D16").Value = "83.73"
$ws.Range("E16").Value = "  +0.81%  "
$ws.Range("D17").Value = "2.098.52"
$ws.Range("E17").Value = "  -0.48%  "
$ws.Range("D18").Value = "6.154"
$ws.Range("E18").Value = "  +0.50%  "
$ws.Range("D19").Value = "29.458.39"
$ws.Range("E19").Value = "  +0.09%  "
$ws.Range("D20").Value = "229.43"
$ws.Range("E20").Value = "  +0.53%  "
$ws.Range("D21").Value = "12.47"
$ws.Range("E21").Value = "  -0.03%  "
$ws.Range("E22").Value = "  +0.09%  "
$ws.Range("D23").Value = "7.433"
$ws.Range("E23").Value = "  -1.38%  "
$ws.Range("E24").Value = "  +0.00%  "
$ws.Range("D25").Value = "157.14"
$ws.Range("D26").Value = "0.1388"
$ws.Range("E26").Value = "  -0.28%  "
$ws.Range("E27").Value = "  +0.24%  "
$ws.Range("D28").Value = "17.71"
$ws.Range("E28").Value = "  +0.26%  "
$ws.Range("E29").Value = "  +4.03%  "
$ws.Range("E30").Value = "  -0.28%  "
$ws.Range("D31").Value = "0.05708"
$ws.Range("E31").Value = "  -0.40%  "
$ws.Range("D32").Value = "4.129"
$ws.Range("E32").Value = "  +0.17%  "
$ws.Range("D33").Value = "4.054"
$ws.Range("E33").Value = "  +0.80%  "
$ws.Range("E34").Value = "  +0.23%  "
$ws.Range("E35").Value = "  +0.75%  "
$ws.Range("D36").Value = "0.7085"
$ws.Range("E36").Value = "  -1.11%  "
$ws.Range("D37").Value = "2.585"
$ws.Range("E37").Value = "  -0.11%  "
$ws.Range("D38").Value = "2.776"
$ws.Range("E38").Value = "  -0.06%  "
$ws.Range("E39").Value = "  -0.78%  "
$ws.Range("D40").Value = "1.220.84"
$ws.Range("E40").Value = "  -2.38%  "
$ws.Range("D41").Value = "6.491"
$ws.Range("E41").Value = "  +5.11%  "
$ws.Range("D42").Value = "0.9107"
$ws.Range("E42").Value = "  +0.20%  "
$ws.Range("E43").Value = "  +0.02%  "
$ws.Range("D44").Value = "2.007.56"
$ws.Range("E44").Value = "  -0.50%  "
$ws.Range("D45").Value = "101.77"
$ws.Range("E45").Value = "  +0.16%  "
$ws.Range("D46").Value = "66.38"
$ws.Range("E46").Value = "  +0.46%  "
$ws.Range("E47").Value = "  +1.74%  "
$ws.Range("D48").Value = "7.137"
$ws.Range("E48").Value = "  +0.53%  "
$ws.Range("E49").Value = "  +0.06%  "
$ws.Range("D50").Value = "8.976"
$ws.Range("E50").Value = "  -1.25%  "
$ws.Range("D51").Value = "1.682"
$ws.Range("E51").Value = "  +0.25%  "
